$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking
# values like "0.9998" or "1.000" keep their exact formatting
# instead of being auto-converted into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.560.96'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '1.935.10'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '246.30'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '0.4851'
$ws.Range("E7").Value = '  +2.60%  '
$ws.Range("D8").Value = '0.2922'
$ws.Range("D9").Value = '0.06805'
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("D10").Value = '112.69'
$ws.Range("E10").Value = '  +6.17%  '
$ws.Range("D11").Value = '19.48'
$ws.Range("D12").Value = '1.937.86'
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = '5.498'
$ws.Range("E13").Value = '  +2.94%  '
$ws.Range("D14").Value = '0.07592'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").Value = '0.6812'
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("D16").Value = '299.68'
$ws.Range("E16").Value = '  +3.56%  '
$ws.Range("D17").Value = '30.567.73'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").Value = '13.12'
$ws.Range("E18").Value = '  +1.14%  '
$ws.Range("D19").Value = '0.000007664'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '5.566'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '0.9990'
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("D22").Value = '2.187.73'
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '6.524'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").Value = '9.561'
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("D26").Value = '168.13'
$ws.Range("E26").Value = '  +0.69%  '
$ws.Range("D27").Value = '20.39'
$ws.Range("E27").Value = '  -1.93%  '
$ws.Range("D28").Value = '2.132'
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D29").Value = '0.1070'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '1.435'
$ws.Range("E30").Value = '  +2.07%  '
$ws.Range("D31").Value = '4.178'
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").Value = '4.109'
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("D33").Value = '0.05012'
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("D34").Value = '0.7495'
$ws.Range("E34").Value = '  +2.04%  '
$ws.Range("D35").Value = '1.149'
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("D36").Value = '0.02037'
$ws.Range("E36").Value = '  -0.94%  '
$ws.Range("D37").Value = '2.714'
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").Value = '2.696'
$ws.Range("E38").Value = '  +0.46%  '
$ws.Range("D39").Value = '2.037'
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("D40").Value = '110.19'
$ws.Range("E40").Value = '  -1.41%  '
$ws.Range("D41").Value = '0.4468'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").Value = '0.8725'
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("D43").Value = '5.840'
$ws.Range("E43").Value = '  -1.19%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '69.80'
$ws.Range("E45").Value = '  +2.64%  '
$ws.Range("D46").Value = '7.300'
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = '49.15'
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").Value = '9.316'
$ws.Range("E48").Value = '  -1.45%  '
$ws.Range("D49").Value = '0.1235'
$ws.Range("E49").Value = '  -2.21%  '
$ws.Range("D50").Value = '0.2548'
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("D51").Value = '35.10'
$ws.Range("E51").Value = '  -0.56%  '
